# Update TPM-derived NATMI LR-pair metrics for Efna1-Epha1
# Sheet layout: rows 2-17 = 4x4 grid of (Sending cluster x Target cluster)
# in the fixed cluster order ECs, FAPs, MuSCs, Resolving-Mac.
# Columns:
#   G = Ligand average expression value   (depends only on Sending cluster)
#   H = Ligand total expression value     (depends only on Sending cluster)
#   I = Ligand derived specificity of avg expr value  = G / sum(G over 4 sending clusters)
#   J = Ligand derived specificity of total expr value = H / sum(H over 4 sending clusters)
#   M = Receptor average expression value (depends only on Target cluster)
#   N = Receptor total expression value   (depends only on Target cluster)
#   O = Receptor derived specificity of avg expr value  = M / sum(M over 4 target clusters)
#   P = Receptor derived specificity of total expr value = N / sum(N over 4 target clusters)
#   Q = Edge average expression weight        = G * M
#   R = Edge total expression weight          = H * N
#   S = Edge average expression derived specificity = Q / sum(Q over all 16 rows)
#   T = Edge total expression derived specificity   = R / sum(R over all 16 rows)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$clusters = @("ECs", "FAPs", "MuSCs", "Resolving-Mac")

# New ligand (sending-cluster) values, keyed by cluster order above
$G = @{
    "ECs"           = 23.73148533333334
    "FAPs"          = 2.471212
    "MuSCs"         = 1.148663
    "Resolving-Mac" = 0.07412966666666666
}
$H = @{
    "ECs"           = 71.194456
    "FAPs"          = 7.413636
    "MuSCs"         = 3.445989
    "Resolving-Mac" = 0.222389
}

# New receptor (target-cluster) values, keyed by cluster order above
$M = @{
    "ECs"           = 5.844648666666667
    "FAPs"          = 4.370261666666667
    "MuSCs"         = 6.452372666666666
    "Resolving-Mac" = 1.570781
}
$N = @{
    "ECs"           = 17.533946
    "FAPs"          = 13.110785
    "MuSCs"         = 19.357118
    "Resolving-Mac" = 4.712343
}

$sumG = 0.0; $sumH = 0.0; $sumM = 0.0; $sumN = 0.0
foreach ($c in $clusters) {
    $sumG += $G[$c]
    $sumH += $H[$c]
    $sumM += $M[$c]
    $sumN += $N[$c]
}

# First pass: compute Q/R per row so we can get the grand totals for S/T
$rows = New-Object System.Collections.ArrayList
$sumQ = 0.0; $sumR = 0.0
$r = 2
foreach ($s in $clusters) {
    foreach ($t in $clusters) {
        $q = $G[$s] * $M[$t]
        $rw = $H[$s] * $N[$t]
        $sumQ += $q
        $sumR += $rw
        $null = $rows.Add(@{ row = $r; s = $s; t = $t; q = $q; rw = $rw })
        $r += 1
    }
}

# Second pass: write every cell
foreach ($entry in $rows) {
    $row = $entry.row
    $s = $entry.s
    $t = $entry.t

    $gVal = $G[$s]
    $hVal = $H[$s]
    $iVal = $gVal / $sumG
    $jVal = $hVal / $sumH

    $mVal = $M[$t]
    $nVal = $N[$t]
    $oVal = $mVal / $sumM
    $pVal = $nVal / $sumN

    $qVal = $entry.q
    $rVal = $entry.rw
    $sVal = $qVal / $sumQ
    $tVal = $rVal / $sumR

    $ws.Range("G$row").Value = $gVal
    $ws.Range("H$row").Value = $hVal
    $ws.Range("I$row").Value = $iVal
    $ws.Range("J$row").Value = $jVal

    $ws.Range("M$row").Value = $mVal
    $ws.Range("N$row").Value = $nVal
    $ws.Range("O$row").Value = $oVal
    $ws.Range("P$row").Value = $pVal

    $ws.Range("Q$row").Value = $qVal
    $ws.Range("R$row").Value = $rVal
    $ws.Range("S$row").Value = $sVal
    $ws.Range("T$row").Value = $tVal
}
